$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H11").Value = 33.5
$ws.Range("I11").Value = 33.5
$ws.Range("K11").Value = 33.5
$ws.Range("M11").Value = 106.5

$ws.Range("H17").Value = 1170.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1170.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3510.6
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").Value = -3846.6

$ws.Range("H51").Value = 12399.5
$ws.Range("J51").Value = 12399
$ws.Range("L51").Value = 12399
$ws.Range("N51").Value = -13367

$ws.Range("H70").Value = 1500
$ws.Range("J70").Value = 1525
$ws.Range("L70").Value = 4575
$ws.Range("N70").Value = -5115

$ws.Range("H73").Value = 1500
$ws.Range("J73").Value = 1525
$ws.Range("L73").Value = 4575
$ws.Range("N73").Value = -6447

$ws.Range("H116").Value = 15750
$ws.Range("J116").Value = 15750
$ws.Range("L116").Value = 15750
$ws.Range("N116").Value = -22634

$ws.Range("H127").Value = 472.5
$ws.Range("I127").Value = 472.5
$ws.Range("K127").Value = 1417.5
$ws.Range("M127").Value = 3542.5

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H5").Value = 551
$ws.Range("J5").Value = 551
$ws.Range("L5").Value = 551
$ws.Range("N5").Value = -775

$ws.Range("H45").Value = 1762.5
$ws.Range("I45").Value = 1762.5
$ws.Range("K45").Value = 1762.5
$ws.Range("M45").Value = -1385.5

$ws.Range("H74").Value = 1080
$ws.Range("I74").Value = 1066.6666
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1066.6666
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -192.6666
$ws.Range("N74").Value = -2848

$ws.Range("H77").Value = 1080
$ws.Range("I77").Value = 1066.6666
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 5333.333000000001
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -965.3330000000005
$ws.Range("N77").Value = -14236

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H4").Value = 551
$ws.Range("J4").Value = 551
$ws.Range("L4").Value = 551
$ws.Range("N4").Value = -781

$ws.Range("H102").Value = 8726.223
$ws.Range("I102").Value = 8726.223
$ws.Range("K102").Value = 8726.223
$ws.Range("M102").Value = -5481.223

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H22").Value = 505.75
$ws.Range("I22").Value = 517.6
$ws.Range("K22").Value = 517.6
$ws.Range("M22").Value = -167.6

$ws.Range("H41").Value = 2625
$ws.Range("I41").Value = 2625
$ws.Range("K41").Value = 2625
$ws.Range("M41").Value = -2197

$ws.Range("H47").Value = 22250
$ws.Range("I47").Value = 19500
$ws.Range("K47").Value = 19500
$ws.Range("M47").Value = -18934

$ws.Range("H62").Value = 7995
$ws.Range("I62").Value = 7995
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7995
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -7371
$ws.Range("N62").ClearContents() | Out-Null

$ws.Range("H65").Value = 7995
$ws.Range("I65").Value = 7995
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 39975
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -36855
$ws.Range("N65").ClearContents() | Out-Null

$ws.Range("H69").Value = 6750
$ws.Range("I69").Value = 6750
$ws.Range("K69").Value = 6750
$ws.Range("M69").Value = -6001

$ws.Range("H72").Value = 6750
$ws.Range("I72").Value = 6750
$ws.Range("K72").Value = 20250
$ws.Range("M72").Value = -16506

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents() | Out-Null

$ws.Range("H122").Value = 2169.2
$ws.Range("I122").Value = 949.6667
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 2849.0001
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = -399.0001000000002
$ws.Range("N122").Value = -16895.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H2").Value = 139.5
$ws.Range("J2").Value = 108.5
$ws.Range("L2").Value = 651
$ws.Range("N2").Value = -877

$ws.Range("H108").Value = 560.8570999999999
$ws.Range("I108").Value = 560.8570999999999
$ws.Range("K108").Value = 1682.5713
$ws.Range("M108").Value = 1197.4287

$ws.Range("H113").Value = 3887.6
$ws.Range("I113").Value = 874.5
$ws.Range("J113").Value = 5896.3335
$ws.Range("K113").Value = 2623.5
$ws.Range("L113").Value = 17689.0005
$ws.Range("M113").Value = -453.5
$ws.Range("N113").Value = -22029.0005

$ws.Range("H117").Value = 7622.1665
$ws.Range("J117").Value = 11086.25
$ws.Range("L117").Value = 33258.75
$ws.Range("N117").Value = -40142.75

$ws.Range("H131").Value = 2090.9092
$ws.Range("J131").Value = 2064.516
$ws.Range("L131").Value = 6193.548000000001
$ws.Range("N131").Value = -16273.548

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 40
$ws.Range("I2").Value = 10
$ws.Range("K2").Value = 10
$ws.Range("M2").Value = 103

$ws.Range("H43").Value = 19135.715
$ws.Range("I43").Value = 4012.75
$ws.Range("J43").Value = 39299.668
$ws.Range("K43").Value = 4012.75
$ws.Range("L43").Value = 39299.668
$ws.Range("M43").Value = -3861.75
$ws.Range("N43").Value = -39601.668

$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19180

$ws.Range("H102").Value = 524.2857
$ws.Range("I102").Value = 414.2
$ws.Range("J102").Value = 799.5
$ws.Range("K102").Value = 414.2
$ws.Range("L102").Value = 799.5
$ws.Range("M102").Value = 1207.8
$ws.Range("N102").Value = -4043.5

$ws.Range("H105").Value = 17886
$ws.Range("J105").Value = 17886
$ws.Range("L105").Value = 17886
$ws.Range("N105").Value = -24874

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents() | Out-Null

$ws.Range("H141").Value = 80429
$ws.Range("J141").Value = 80429
$ws.Range("L141").Value = 80429
$ws.Range("N141").Value = -90789

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("M22").ClearContents() | Out-Null
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -637
$ws.Range("M27").ClearContents() | Out-Null
$ws.Range("N27").Value = -1214

$ws.Range("H40").Value = 3750
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -5272

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents() | Out-Null

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents() | Out-Null

$ws.Range("H122").Value = 15000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents() | Out-Null

$ws.Range("H132").Value = 5707.769
$ws.Range("I132").Value = 4899.1665
$ws.Range("K132").Value = 14697.4995
$ws.Range("M132").Value = -12167.4995

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H4").Value = 18519.166
$ws.Range("J4").Value = 18519.166
$ws.Range("L4").Value = 18519.166
$ws.Range("N4").Value = -18745.166

$ws.Range("H58").Value = 32716.6
$ws.Range("I58").Value = 27872.25
$ws.Range("K58").Value = 27872.25
$ws.Range("M58").Value = -27564.25

$ws.Range("H122").Value = 21099.2
$ws.Range("I122").Value = 1749
$ws.Range("K122").Value = 5247
$ws.Range("M122").Value = -2797

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents() | Out-Null

